# Apply the "Add files via upload" edit to training.xlsx
#
# Summary of the change:
#  - Four row labels in column A are renamed from "...職員總數" to "...職員人數"
#    (rows 7, 11, 15, 19).
#  - A new row (26) is appended, re-using the existing "主管預算數" figure in
#    column B but labelled "113年工務局主管決算數" in column A.
#  - The B28 cell becomes the active selection afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the four "...職員總數" labels to "...職員人數" ---------------
$ws.Range("A7").Value  = "113年新建工程處職員人數"
$ws.Range("A11").Value = "113年道路養護工程處職員人數"
$ws.Range("A15").Value = "113年公園處職員人數"
$ws.Range("A19").Value = "113年違章建築處理大隊職員人數"

# --- 2. Append a new row 26 (主管決算數), copying the style of row 25 ------
$budgetText = $ws.Range("B23").Value()

$ws.Range("A25:B25").Copy()
$ws.Range("A26:B26").PasteSpecial(-4122)   # xlPasteFormats (style + number format only)

$ws.Range("A26").Value = "113年工務局主管決算數"
$ws.Range("B26").Value = $budgetText
$ws.Rows.Item(26).RowHeight = $ws.Rows.Item(25).RowHeight()

# --- 3. Update the current selection to B28 (matches the saved view) ------
$ws.Activate()
$ws.Range("B28").Select()
